# Update gh-pages to output generated at 456a3b4
# Applies the numeric/text refreshes captured in the upstream diff to the
# "展览" and "全部类型" worksheets of the workbook.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---------------------------------------------------------
$wsExhibit.Range("F2").Value = 132
$wsExhibit.Range("F3").Value = 128
$wsExhibit.Range("F7").Value = 975
$wsExhibit.Range("F8").Value = 938
$wsExhibit.Range("F10").Value = 96
$wsExhibit.Range("F13").Value = 922
$wsExhibit.Range("F14").Value = 1802
$wsExhibit.Range("F15").Value = 3931
$wsExhibit.Range("F16").Value = 1161
$wsExhibit.Range("F17").Value = 111
$wsExhibit.Range("F18").Value = 2594
$wsExhibit.Range("F20").Value = 1078
$wsExhibit.Range("F21").Value = 3575
$wsExhibit.Range("F22").Value = 754
$wsExhibit.Range("F23").Value = 839
$wsExhibit.Range("F25").Value = 2214
$wsExhibit.Range("F26").Value = 111

$wsExhibit.Range("D27").Value = "丁城路丁桥桃花湖公园北区 典酷沉浸式艺术空间(桃花湖店)"
$wsExhibit.Range("F27").Value = 837

$wsExhibit.Range("F28").Value = 168
$wsExhibit.Range("F29").Value = 475
$wsExhibit.Range("F30").Value = 200

$wsExhibit.Range("C31").Value = "杭州·造梦探险家——二次元同好会（取消）"
$wsExhibit.Range("G31").Value = "不可售"

$wsExhibit.Range("F32").Value = 1337
$wsExhibit.Range("F33").Value = 1949
$wsExhibit.Range("F35").Value = 43
$wsExhibit.Range("F38").Value = 277
$wsExhibit.Range("F39").Value = 40

# --- 全部类型 (sheet4) ------------------------------------------------------
$wsAll.Range("F2").Value = 430
$wsAll.Range("F3").Value = 132
$wsAll.Range("F6").Value = 975
$wsAll.Range("F7").Value = 938
$wsAll.Range("F10").Value = 96
$wsAll.Range("F14").Value = 922
$wsAll.Range("F15").Value = 1802
$wsAll.Range("F16").Value = 3931
$wsAll.Range("F17").Value = 1161
$wsAll.Range("F18").Value = 111
$wsAll.Range("F20").Value = 2594
$wsAll.Range("F21").Value = 1078
$wsAll.Range("F22").Value = 3575
$wsAll.Range("F23").Value = 754
$wsAll.Range("F24").Value = 839
$wsAll.Range("F27").Value = 2214
$wsAll.Range("F31").Value = 111

$wsAll.Range("D33").Value = "丁城路丁桥桃花湖公园北区 典酷沉浸式艺术空间(桃花湖店)"
$wsAll.Range("F33").Value = 837

$wsAll.Range("F34").Value = 168
$wsAll.Range("F35").Value = 475
$wsAll.Range("F36").Value = 200

# Row 37 is fully replaced with a new event (old "造梦探险家" entry is gone
# from this sheet; it became the "英雄时代2024哈瓦西钢琴演奏会" event).
# B37 must stay a plain text date (like the rest of column B), so force a
# text number format before assigning, then restore the original style, to
# keep Excel's auto date-detection from turning the string into a serial.
$dateCell = $wsAll.Range("B37")
$origStyle = $dateCell.Style
$dateCell.NumberFormat = "@"
$dateCell.Value = "2024-06-05"
$dateCell.Style = $origStyle

$wsAll.Range("C37").Value = "杭州·英雄时代2024哈瓦西钢琴演奏会"
$wsAll.Range("D37").Value = "中国杭州北山路86号西湖岳湖景区 中国杭州西湖岳湖景区印象西湖"
$wsAll.Range("E37").Value = "2024.06.05 20:00-06.05 21:30"
$wsAll.Range("F37").Value = 1
$wsAll.Range("G37").Value = 499
$wsAll.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=83902"
$wsAll.Range("I37").Value = "//i2.hdslb.com/bfs/openplatform/202404/BFRFmKpT1712569969076.jpeg"

$wsAll.Range("F38").Value = 1337
$wsAll.Range("F39").Value = 1949
$wsAll.Range("F43").Value = 43
$wsAll.Range("F45").Value = 277
$wsAll.Range("F46").Value = 40
